$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
# Swap "to" (D1) and "moved_out" (E1)
$ws.Range("D1").Value = "moved_out"
$ws.Range("E1").Value = "to"
# Rename nij_score -> nij
$ws.Range("H1").Value = "nij"

# --- Row 2 ---
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = "二食堂"
$ws.Range("H2").Value = 0.4283840080205571

# --- Row 3 ---
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = "二食堂"
$ws.Range("F3").Value = 10
$ws.Range("G3").Value = 10
$ws.Range("H3").Value = 0.5017841360825993

# --- Row 4 ---
$ws.Range("C4").Value = "二食堂"
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = "一食堂"
$ws.Range("F4").Value = 10
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0.3948432400342524

# --- Row 5 ---
$ws.Range("C5").Value = "东门"
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = "三食堂"
$ws.Range("F5").Value = 20
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 0.4205826887426308
